# LOM3046.docx content-rotation edit
# The paragraph/run/style skeleton of the document does not change; only the
# text carried by certain runs is swapped around (a rotation of several
# blocks of text among paragraphs), plus the "Critério:"/"Norma de
# recuperação:" run bodies inside the "Avaliação" bullet paragraph are
# rebuilt to host the old Bibliografia list (with the old exam-recovery text
# moving down to where the bibliography list used to start).
#
# Every Find.Execute() below is scoped to a single paragraph's Range (or an
# even narrower sub-range within that paragraph) so that identical text that
# is "passing through" several paragraphs at once during this rotation can
# never be matched ambiguously.

$d = $word.ActiveDocument

function Replace-InRange($range, [string]$old, [string]$new) {
    $ok = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replace failed: old=[$old] not found"
    }
}

# --- 1. "Objetivos" body paragraph (index 6) ---------------------------
# old: original "Objetivos" text -> new: old "Programa resumido" text
Replace-InRange $d.Paragraphs.Item(6).Range `
    "Apresentação introdutória das técnicas de análise microestrutural de materiais. Apresentação das técnicas e equipamentos necessários para a análise microestrutural. Seleção adequada das técnicas experimentais. Apresentação das técnicas adequadas de preparação de amostras. Verificação dos custos envolvidos nas técnicas de caracterização microestrutural." `
    "A Microestrutura dos Materiais. Difratometria de raios X. Análise Microestrutural utilizando Luz Síncrotron. Microscopia Óptica. Microscopia Eletrônica. Microscopia de Tunelamento e de Força Atômica. Análise Química de Microrregiões. Análises Térmicas. Fluorescência de raios X. Técnicas Indiretas de Análise de Microestrutura. Seleção de Técnicas Experimentais."

# --- 2. "Docente(s) Responsável(eis)" list paragraph (index 8) ----------
# run 1: old docente name -> old "Objetivos" text (keeps its trailing break)
Replace-InRange $d.Paragraphs.Item(8).Range `
    "6495737 - Durval Rodrigues Junior" `
    "Apresentação introdutória das técnicas de análise microestrutural de materiais. Apresentação das técnicas e equipamentos necessários para a análise microestrutural. Seleção adequada das técnicas experimentais. Apresentação das técnicas adequadas de preparação de amostras. Verificação dos custos envolvidos nas técnicas de caracterização microestrutural."

# run 2: old docente name -> old "Programa" text
Replace-InRange $d.Paragraphs.Item(8).Range `
    "1643715 - Paulo Atsushi Suzuki" `
    "1. A Microestrutura dos Materiais. 2. Difratometria de raios X. 3. Análise Microestrutural utilizando Luz Síncrotron. 4. Microscopia Óptica. 5. Microscopia Eletrônica. 6. Microscopia de Tunelamento e de Força Atômica. 7. Análise Química de Microrregiões. 8. Análises Térmicas. 9. Fluorescência de raios X. 10. Técnicas Indiretas de Análise de Microestrutura. 11. Seleção de Técnicas Experimentais."

# --- 3. "Programa resumido" body paragraph (index 10) -------------------
Replace-InRange $d.Paragraphs.Item(10).Range `
    "A Microestrutura dos Materiais. Difratometria de raios X. Análise Microestrutural utilizando Luz Síncrotron. Microscopia Óptica. Microscopia Eletrônica. Microscopia de Tunelamento e de Força Atômica. Análise Química de Microrregiões. Análises Térmicas. Fluorescência de raios X. Técnicas Indiretas de Análise de Microestrutura. Seleção de Técnicas Experimentais." `
    "Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre)."

# --- 4. "Programa" body paragraph (index 12) -----------------------------
Replace-InRange $d.Paragraphs.Item(12).Range `
    "1. A Microestrutura dos Materiais. 2. Difratometria de raios X. 3. Análise Microestrutural utilizando Luz Síncrotron. 4. Microscopia Óptica. 5. Microscopia Eletrônica. 6. Microscopia de Tunelamento e de Força Atômica. 7. Análise Química de Microrregiões. 8. Análises Térmicas. 9. Fluorescência de raios X. 10. Técnicas Indiretas de Análise de Microestrutura. 11. Seleção de Técnicas Experimentais." `
    "A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2."

# --- 5. "Avaliação" bullet paragraph (index 14) -------------------------
# Work on the "Norma de recuperação:" body run FIRST: its current text is the
# same string that will be written into the "Método:" body run later, so it
# must be swapped out before that happens (otherwise the later Find could
# match the freshly-written "Método:" text instead of this one).
Replace-InRange $d.Paragraphs.Item(14).Range `
    "Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação)." `
    "6495737 - Durval Rodrigues Junior"

# "Critério:" body run: old "A Nota Final..." sentence -> the 19-entry
# bibliography list (each entry separated by a manual line break, i.e. a
# <w:br/>), with a trailing manual line break before the (unchanged) "Norma
# de recuperação: " bold label that follows it.
$biblioReplacement = "1. Van Vlack, L.H. Princípios de Ciência e Tecnologia dos Materiais, 4a.ed., Ed. Campus, Rio de Janeiro, 1984. ^l2. Shackelford, J.F. Introduction to Materials Science for Engineers. 4th Edition. Prentice Hall Inc., 1996. ^l3. Padilha, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985. ^l4. Guy, A.G. Ciência dos Materiais. Livros Técnicos e Científicos Editora, 1982. ^l5. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. ^l6. Nondestructive Characterization of Materials. Series. Plenum Press, New York. ^l7. Yacobi, B.G. Holt, D.B. Kazmerski, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994. ^l8. Lowell, S.; Shields, J. E.; Thomas, M. A.; Thommes, M. Characterization of Porous Solids and Powders: Surface Area, Pore Size and Density, Springer, 2010. ^l9. Murphy, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001. ^l10. Wu, Q.; Merchant, F.; Castleman, K. Microscope Image Processing, Academic Press, 2008. ^l11. Cullity, B. D.; Stock, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001. ^l12. Goldstein, J.; et al., Scanning Electron Microscopy and X-ray Microanalysis, Springer, 2003. ^l13. Hatakeyama, T.; Zhenhai, L. Handbook of Thermal Analysis, NY: Wiley, 1999. ^l14. Haines, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002. ^l15. Schramm, G. Reologia e Reometria. Editora Artliber, 2006.^l16. Azevedo, A. D.; Mothe, C. G. Análise Térmica de Materiais. São Paulo: ARTLIBER, 2009.^l17. Brown, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.^l18. Muller, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.^l19. Speyer, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.^l"
Replace-InRange $d.Paragraphs.Item(14).Range `
    "A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2." `
    $biblioReplacement

# "Método:" body run: old exam text -> old "Norma de recuperação" exam text
Replace-InRange $d.Paragraphs.Item(14).Range `
    "Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre)." `
    "Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação)."

# --- 6. "Bibliografia" body paragraph (index 16) ------------------------
$biblioSearch = "1. Van Vlack, L.H. Princípios de Ciência e Tecnologia dos Materiais, 4a.ed., Ed. Campus, Rio de Janeiro, 1984. ^l2. Shackelford, J.F. Introduction to Materials Science for Engineers. 4th Edition. Prentice Hall Inc., 1996. ^l3. Padilha, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985. ^l4. Guy, A.G. Ciência dos Materiais. Livros Técnicos e Científicos Editora, 1982. ^l5. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. ^l6. Nondestructive Characterization of Materials. Series. Plenum Press, New York. ^l7. Yacobi, B.G. Holt, D.B. Kazmerski, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994. ^l8. Lowell, S.; Shields, J. E.; Thomas, M. A.; Thommes, M. Characterization of Porous Solids and Powders: Surface Area, Pore Size and Density, Springer, 2010. ^l9. Murphy, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001. ^l10. Wu, Q.; Merchant, F.; Castleman, K. Microscope Image Processing, Academic Press, 2008. ^l11. Cullity, B. D.; Stock, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001. ^l12. Goldstein, J.; et al., Scanning Electron Microscopy and X-ray Microanalysis, Springer, 2003. ^l13. Hatakeyama, T.; Zhenhai, L. Handbook of Thermal Analysis, NY: Wiley, 1999. ^l14. Haines, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002. ^l15. Schramm, G. Reologia e Reometria. Editora Artliber, 2006.^l16. Azevedo, A. D.; Mothe, C. G. Análise Térmica de Materiais. São Paulo: ARTLIBER, 2009.^l17. Brown, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.^l18. Muller, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.^l19. Speyer, R. Thermal analysis of materials, New York: Marcel Dekker, 1994."
Replace-InRange $d.Paragraphs.Item(16).Range `
    $biblioSearch `
    "1643715 - Paulo Atsushi Suzuki"

Write-Output "done"
